$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $origStyle
}

Set-TextCell $ws "D2" "65.681.27"
Set-TextCell $ws "E2" "  -0.83%  "
Set-TextCell $ws "D3" "3.451.38"
Set-TextCell $ws "E3" "  -3.11%  "
Set-TextCell $ws "D5" "594.97"
Set-TextCell $ws "E5" "  -1.58%  "
Set-TextCell $ws "D6" "136.75"
Set-TextCell $ws "E6" "  -7.02%  "
Set-TextCell $ws "D7" "3.448.57"
Set-TextCell $ws "E7" "  -3.19%  "
Set-TextCell $ws "E8" "  -0.12%  "
Set-TextCell $ws "E9" "  +0.59%  "
Set-TextCell $ws "D10" "7.44"
Set-TextCell $ws "E10" "  -5.34%  "
Set-TextCell $ws "E11" "  -8.94%  "
Set-TextCell $ws "E12" "  -7.35%  "
Set-TextCell $ws "D13" "4.031.31"
Set-TextCell $ws "E13" "  -3.44%  "
Set-TextCell $ws "D14" "0.0000182"
Set-TextCell $ws "E14" "  -10.15%  "
Set-TextCell $ws "D15" "26.73"
Set-TextCell $ws "E15" "  -8.66%  "
Set-TextCell $ws "D16" "3.460.99"
Set-TextCell $ws "E16" "  -2.54%  "
Set-TextCell $ws "D17" "65.593.15"
Set-TextCell $ws "E17" "  -1.01%  "
Set-TextCell $ws "E18" "  -2.04%  "
Set-TextCell $ws "D19" "9.90"
Set-TextCell $ws "E19" "  -9.82%  "
Set-TextCell $ws "D20" "5.79"
Set-TextCell $ws "E20" "  -7.74%  "
Set-TextCell $ws "D21" "13.71"
Set-TextCell $ws "E21" "  -6.66%  "
Set-TextCell $ws "D22" "394.09"
Set-TextCell $ws "E22" "  -6.22%  "
Set-TextCell $ws "E23" "  -9.47%  "
Set-TextCell $ws "D24" "73.33"
Set-TextCell $ws "E24" "  -5.87%  "
Set-TextCell $ws "E25" "  -0.01%  "
Set-TextCell $ws "D26" "3.591.25"
Set-TextCell $ws "E26" "  -3.30%  "
Set-TextCell $ws "E27" "  -9.34%  "
Set-TextCell $ws "D28" "1.01"
Set-TextCell $ws "E28" "  +0.58%  "
Set-TextCell $ws "D29" "2.27"
Set-TextCell $ws "E29" "  -8.38%  "
Set-TextCell $ws "D30" "7.25"
Set-TextCell $ws "E30" "  -9.17%  "
Set-TextCell $ws "D31" "8.23"
Set-TextCell $ws "E31" "  -11.27%  "
Set-TextCell $ws "D32" "3.454.02"
Set-TextCell $ws "E32" "  -3.16%  "
Set-TextCell $ws "E33" "  -0.01%  "
Set-TextCell $ws "D34" "0.147"
Set-TextCell $ws "E34" "  -6.08%  "
Set-TextCell $ws "D35" "22.99"
Set-TextCell $ws "E35" "  -6.94%  "
Set-TextCell $ws "D36" "172.35"
Set-TextCell $ws "E36" "  -1.68%  "
Set-TextCell $ws "E37" "  -12.14%  "
Set-TextCell $ws "E38" "  -9.80%  "
Set-TextCell $ws "D39" "1.52"
Set-TextCell $ws "E39" "  -6.58%  "
Set-TextCell $ws "D40" "4.81"
Set-TextCell $ws "E40" "  -11.24%  "
Set-TextCell $ws "D41" "0.0776"
Set-TextCell $ws "E41" "  -7.65%  "
Set-TextCell $ws "D42" "0.824"
Set-TextCell $ws "E42" "  -5.66%  "
Set-TextCell $ws "D43" "43.49"
Set-TextCell $ws "E43" "  -5.20%  "
Set-TextCell $ws "E44" "  -0.05%  "
Set-TextCell $ws "D45" "4.43"
Set-TextCell $ws "D46" "1.62"
Set-TextCell $ws "E46" "  -11.27%  "
Set-TextCell $ws "E47" "  -1.53%  "
Set-TextCell $ws "D48" "1.11"
Set-TextCell $ws "E48" "  -1.61%  "
Set-TextCell $ws "D49" "6.59"
Set-TextCell $ws "E49" "  -7.06%  "
Set-TextCell $ws "D50" "2.11"
Set-TextCell $ws "E50" "  -14.04%  "
Set-TextCell $ws "D51" "2.208.19"
Set-TextCell $ws "E51" "  -7.32%  "
